$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Range("H33").Value = 285.77777
$ws.Range("I33").Value = 131.57143
$ws.Range("K33").Value = 131.57143
$ws.Range("M33").Value = 97.42857000000001

$ws.Range("H100").Value = 1233.3334
$ws.Range("I100").Value = 1233.3334
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1233.3334
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -692.3334
$ws.Range("N100").ClearContents()

$ws.Range("H111").Value = 886.1818
$ws.Range("I111").Value = 874.8
$ws.Range("J111").Value = 1000
$ws.Range("K111").Value = 2624.4
$ws.Range("L111").Value = 3000
$ws.Range("M111").Value = 442.6000000000004
$ws.Range("N111").Value = -9134

$ws.Range("H120").Value = 24944.445
$ws.Range("J120").Value = 24944.445
$ws.Range("L120").Value = 24944.445
$ws.Range("N120").Value = -34620.445

$ws.Range("H129").Value = 806235.3
$ws.Range("J129").Value = 950866.8
$ws.Range("L129").Value = 2852600.4
$ws.Range("N129").Value = -2862600.4

$ws.Range("H137").Value = 768.6739
$ws.Range("I137").Value = 687.5946
$ws.Range("K137").Value = 2062.7838
$ws.Range("M137").Value = 487.2161999999998

$ws.Range("H138").Value = 2540.74
$ws.Range("I138").Value = 1070.75
$ws.Range("J138").Value = 3520.7334
$ws.Range("K138").Value = 3212.25
$ws.Range("L138").Value = 10562.2002
$ws.Range("M138").Value = 1927.75
$ws.Range("N138").Value = -20842.2002

$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Range("H20").Value = 27945.852
$ws.Range("I20").Value = 76808.75
$ws.Range("J20").Value = 7372
$ws.Range("K20").Value = 76808.75
$ws.Range("L20").Value = 7372
$ws.Range("M20").Value = -76561.75
$ws.Range("N20").Value = -7866

$ws.Range("H120").Value = 49000
$ws.Range("J120").Value = 49000
$ws.Range("L120").Value = 49000
$ws.Range("N120").Value = -58676

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Range("H41").Value = 16999
$ws.Range("J41").Value = 16999
$ws.Range("L41").Value = 16999
$ws.Range("N41").Value = -17855

$ws.Range("H68").Value = 15519.8
$ws.Range("J68").Value = 15519.8
$ws.Range("L68").Value = 15519.8
$ws.Range("N68").Value = -17017.8

$ws.Range("H71").Value = 15519.8
$ws.Range("J71").Value = 15519.8
$ws.Range("L71").Value = 46559.39999999999
$ws.Range("N71").Value = -54047.39999999999

$ws.Range("H74").Value = 13697
$ws.Range("J74").Value = 13697
$ws.Range("L74").Value = 13697
$ws.Range("N74").Value = -15445

$ws.Range("H77").Value = 13697
$ws.Range("J77").Value = 13697
$ws.Range("L77").Value = 41091
$ws.Range("N77").Value = -49827

$ws.Range("H99").Value = 2832.1538
$ws.Range("I99").Value = 2858.4666
$ws.Range("J99").Value = 2744.4443
$ws.Range("K99").Value = 2858.4666
$ws.Range("L99").Value = 2744.4443
$ws.Range("M99").Value = -1360.4666
$ws.Range("N99").Value = -5740.4443

$ws.Range("H119").Value = 32331.666
$ws.Range("J119").Value = 32331.666
$ws.Range("L119").Value = 32331.666
$ws.Range("N119").Value = -42007.666

$ws.Range("H126").Value = 2832.1538
$ws.Range("I126").Value = 2858.4666
$ws.Range("J126").Value = 2744.4443
$ws.Range("K126").Value = 8575.399800000001
$ws.Range("L126").Value = 8233.332900000001
$ws.Range("M126").Value = -6105.399800000001
$ws.Range("N126").Value = -13173.3329

$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Range("H80").Value = 7268.091
$ws.Range("I80").Value = 7319.6665
$ws.Range("J80").Value = 7248.75
$ws.Range("K80").Value = 21958.9995
$ws.Range("L80").Value = 21746.25
$ws.Range("M80").Value = -21022.9995
$ws.Range("N80").Value = -23618.25

$ws.Range("H83").Value = 7268.091
$ws.Range("I83").Value = 7319.6665
$ws.Range("J83").Value = 7248.75
$ws.Range("K83").Value = 65876.9985
$ws.Range("L83").Value = 65238.75
$ws.Range("M83").Value = -61196.9985
$ws.Range("N83").Value = -74598.75

$ws.Range("H92").Value = 30284.6
$ws.Range("I92").Value = 100234
$ws.Range("J92").Value = 306.2857
$ws.Range("K92").Value = 300702
$ws.Range("L92").Value = 918.8571000000001
$ws.Range("M92").Value = -299454
$ws.Range("N92").Value = -3414.8571

$ws.Range("H121").Value = 822.4737
$ws.Range("J121").Value = 848.7222
$ws.Range("L121").Value = 2546.1666
$ws.Range("N121").Value = -5166.1666

$ws.Range("H131").Value = 17311594
$ws.Range("J131").Value = 47003.816
$ws.Range("L131").Value = 141011.448
$ws.Range("N131").Value = -151091.448

$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Range("H70").Value = 4305.44
$ws.Range("I70").Value = 4217.4707
$ws.Range("K70").Value = 4217.4707
$ws.Range("M70").Value = -3947.4707

$ws.Range("H73").Value = 4305.44
$ws.Range("I73").Value = 4217.4707
$ws.Range("K73").Value = 4217.4707
$ws.Range("M73").Value = -3281.4707

$ws.Range("H119").Value = 26791.2
$ws.Range("J119").Value = 26791.2
$ws.Range("L119").Value = 26791.2
$ws.Range("N119").Value = -36467.2

$ws.Range("H126").Value = 4703
$ws.Range("I126").Value = 4703
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 14109
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11639
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Range("H7").Value = 2473.8386
$ws.Range("I7").Value = 2363.75
$ws.Range("K7").Value = 2363.75
$ws.Range("M7").Value = -2251.75

$ws.Range("H61").Value = 1943.6923
$ws.Range("I61").Value = 1802.3
$ws.Range("J61").Value = 2415
$ws.Range("K61").Value = 1802.3
$ws.Range("L61").Value = 2415
$ws.Range("M61").Value = -1600.3
$ws.Range("N61").Value = -2819

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H113").Value = 1943.6923
$ws.Range("I113").Value = 1802.3
$ws.Range("J113").Value = 2415
$ws.Range("K113").Value = 1802.3
$ws.Range("L113").Value = 2415
$ws.Range("M113").Value = 367.7
$ws.Range("N113").Value = -6755

$ws.Range("H116").Value = 6689993.5
$ws.Range("J116").Value = 6689993.5
$ws.Range("L116").Value = 6689993.5
$ws.Range("N116").Value = -6699171.5

$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800

$ws.Range("H126").Value = 2473.8386
$ws.Range("I126").Value = 2363.75
$ws.Range("K126").Value = 7091.25
$ws.Range("M126").Value = -4621.25

$ws.Range("H127").Value = 47000
$ws.Range("J127").Value = 47000
$ws.Range("L127").Value = 47000
$ws.Range("N127").Value = -56920

$ws.Range("H128").Value = 50000
$ws.Range("J128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("N128").Value = -59960

$ws.Range("H130").Value = 47950
$ws.Range("J130").Value = 47950
$ws.Range("L130").Value = 47950
$ws.Range("N130").Value = -57990

$ws.Range("H132").Value = 2928.255
$ws.Range("I132").Value = 1837.697
$ws.Range("J132").Value = 4927.6113
$ws.Range("K132").Value = 5513.090999999999
$ws.Range("L132").Value = 14782.8339
$ws.Range("M132").Value = -2983.090999999999
$ws.Range("N132").Value = -19842.8339

$ws.Range("H136").Value = 2412.25
$ws.Range("I136").Value = 1462.7384
$ws.Range("K136").Value = 4388.2152
$ws.Range("M136").Value = -1838.2152

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Range("H107").Value = 611.75757
$ws.Range("I107").Value = 577.11536
$ws.Range("J107").Value = 740.4286
$ws.Range("K107").Value = 1731.34608
$ws.Range("L107").Value = 2221.2858
$ws.Range("M107").Value = 188.65392
$ws.Range("N107").Value = -6061.2858

$ws.Range("H116").Value = 10500
$ws.Range("J116").Value = 10500
$ws.Range("L116").Value = 10500
$ws.Range("N116").Value = -19678

$ws.Range("H119").Value = 37641.43
$ws.Range("J119").Value = 37641.43
$ws.Range("L119").Value = 37641.43
$ws.Range("N119").Value = -47317.43
